# Applies the "Light requirement & some state transition flags" edit:
# adds three new columns (M: heat_capacity, N: min_skylight, O: max_skylight)
# with a header row and per-plant values, and updates a couple of cosmetic
# view settings (zoom / selection / column widths).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (match the formatting of the existing K1/L1 header) ---
$ws.Range("K1").Copy($ws.Range("M1"))
$ws.Range("K1").Copy($ws.Range("N1"))
$ws.Range("K1").Copy($ws.Range("O1"))
$ws.Cells.Item(1, 13).Value = "heat_capacity"
$ws.Cells.Item(1, 14).Value = "min_skylight"
$ws.Cells.Item(1, 15).Value = "max_skylight"

# --- New column widths (stored width = ColumnWidth + 5/7) -------------
$ws.Columns.Item(13).ColumnWidth = 13.83203125 - (5/7)
$ws.Columns.Item(14).ColumnWidth = 13.1640625 - (5/7)
$ws.Columns.Item(15).ColumnWidth = 13 - (5/7)

# --- New per-row data: row, heat_capacity, min_skylight, max_skylight -
$plantData = @(
    @(2, 10, 12, 15),
    @(3, 10, 9, 15),
    @(4, 10, 0, 15),
    @(5, 10, 9, 15),
    @(6, 10, 0, 4),
    @(7, 10, 0, 4),
    @(8, 10, 0, 4),
    @(9, 10, 14, 15),
    @(10, 10, 0, 15),
    @(11, 10, 14, 15),
    @(12, 10, 14, 15),
    @(13, 10, 14, 15),
    @(14, 10, 14, 15),
    @(15, 10, 14, 15),
    @(16, 10, 14, 15),
    @(17, 10, 14, 15),
    @(18, 10, 14, 15),
    @(19, 10, 14, 15),
    @(20, 10, 14, 15),
    @(21, 10, 0, 4),
    @(22, 10, 0, 4),
    @(23, 10, 0, 4),
    @(24, 10, 0, 4),
    @(25, 10, 9, 15),
    @(26, 10, 9, 15),
    @(27, 10, 9, 15),
    @(28, 10, 9, 15),
    @(29, 10, 9, 15),
    @(30, 10, 9, 15),
    @(31, 10, 9, 15),
    @(32, 10, 9, 15),
    @(33, 10, 9, 15),
    @(34, 10, 9, 15),
    @(35, 10, 9, 15),
    @(36, 10, 9, 15),
    @(37, 10, 9, 15),
    @(38, 10, 9, 15),
    @(39, 10, 9, 15),
    @(40, 10, 9, 15),
    @(41, 10, 9, 15),
    @(42, 10, 12, 15),
    @(43, 10, 12, 15),
    @(44, 10, 12, 15),
    @(45, 10, 0, 15),
    @(46, 10, 0, 15),
    @(47, 10, 0, 15),
    @(48, 10, 0, 15),
    @(49, 10, 0, 15),
    @(50, 10, 9, 15),
    @(51, 10, 9, 15),
    @(52, 10, 9, 15),
    @(53, 10, 9, 15),
    @(54, 10, 9, 15),
    @(55, 10, 9, 15),
    @(56, 10, 9, 15),
    @(57, 10, 9, 15),
    @(58, 10, 9, 15),
    @(59, 10, 9, 15),
    @(60, 100, 9, 15),
    @(61, 10, 12, 15),
    @(62, 10, 9, 15),
    @(63, 10, 12, 15),
    @(64, 10, 12, 15),
    @(65, 10, 9, 15),
    @(66, 10, 9, 15),
    @(67, 10, 12, 15),
    @(68, 10, 9, 15),
    @(69, 10, 9, 15),
    @(70, 10, 9, 15),
    @(71, 10, 12, 15),
    @(72, 10, 0, 15),
    @(73, 10, 0, 15),
    @(74, 10, 14, 15),
    @(75, 10, 14, 15),
    @(76, 10, 14, 15),
    @(77, 10, 14, 15),
    @(78, 10, 14, 15),
    @(79, 10, 14, 15),
    @(80, 10, 14, 15),
    @(81, 10, 14, 15),
    @(82, 10, 14, 15),
    @(83, 10, 0, 4),
    @(84, 10, 0, 4),
    @(85, 10, 0, 4),
    @(86, 10, 6, 15),
    @(87, 10, 6, 15),
    @(88, 10, 6, 15),
    @(89, 10, 6, 15),
    @(90, 10, 6, 15),
    @(91, 10, 6, 15),
    @(92, 10, 12, 15),
    @(93, 10, 12, 15),
    @(94, 10, 9, 15),
    @(95, 10, 9, 15),
    @(96, 10, 12, 15),
    @(97, 10, 12, 15),
    @(98, 10, 12, 15),
)

foreach ($row in $plantData) {
    $r = $row[0]
    $ws.Cells.Item($r, 13).Value = $row[1]
    $ws.Cells.Item($r, 14).Value = $row[2]
    $ws.Cells.Item($r, 15).Value = $row[3]
}

# --- View state ---------------------------------------------------------
$ws.Range("N64").Select()
$excel.ActiveWindow.Zoom = 124

Write-Host "Applied plant_temperature light-requirement update."
